# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 170-171) above the existing most-recent
# record for this market/product, shifting every subsequent row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 170 (shifts rows 170..218 down to 172..220)
$ws.Rows.Item(170).Resize(2).Insert(-4121)

# New row 170: Mango, "Sin especificar" / "Primera", 2023-03-03 (serial 44988)
$ws.Cells.Item(170, 1).Value = 1
$ws.Cells.Item(170, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(170, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(170, 4).Value = 44988
$ws.Cells.Item(170, 5).Value = 15
$ws.Cells.Item(170, 6).Value = "Fruta"
$ws.Cells.Item(170, 7).Value = 100108
$ws.Cells.Item(170, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(170, 9).Value = 100108002
$ws.Cells.Item(170, 10).Value = "Mango"
$ws.Cells.Item(170, 11).Value = "Sin especificar"
$ws.Cells.Item(170, 12).Value = "Primera"
$ws.Cells.Item(170, 13).Value = 950
$ws.Cells.Item(170, 14).Value = 4500
$ws.Cells.Item(170, 15).Value = 5000
$ws.Cells.Item(170, 16).Value = 4763
$ws.Cells.Item(170, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(170, 18).Value = "Perú"
$ws.Cells.Item(170, 19).Value = 1191
$ws.Cells.Item(170, 20).Value = 4

# New row 171: Mango, "Sin especificar" / "Segunda", 2023-03-03 (serial 44988)
$ws.Cells.Item(171, 1).Value = 1
$ws.Cells.Item(171, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(171, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(171, 4).Value = 44988
$ws.Cells.Item(171, 5).Value = 15
$ws.Cells.Item(171, 6).Value = "Fruta"
$ws.Cells.Item(171, 7).Value = 100108
$ws.Cells.Item(171, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(171, 9).Value = 100108002
$ws.Cells.Item(171, 10).Value = "Mango"
$ws.Cells.Item(171, 11).Value = "Sin especificar"
$ws.Cells.Item(171, 12).Value = "Segunda"
$ws.Cells.Item(171, 13).Value = 830
$ws.Cells.Item(171, 14).Value = 4500
$ws.Cells.Item(171, 15).Value = 5000
$ws.Cells.Item(171, 16).Value = 4711
$ws.Cells.Item(171, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(171, 18).Value = "Perú"
$ws.Cells.Item(171, 19).Value = 1178
$ws.Cells.Item(171, 20).Value = 4
